# TODO.xlsx — "updated TODO for current sprint"
#
# Sprint 1 status updates:
#  - "Mitarbeiter(mit jeweiliger Rolle) anlegen" is now done (100%) and its
#    "Soll" note is replaced by a note about the finished admin feature.
#  - "Zimmer buchen" is now fully done (100%) and its open "Soll" issue note
#    is cleared since the bug was fixed.
#
# Sprint 2 is populated with the 4 new backlog items that were worked on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sprint 1 row 6: "Mitarbeiter(mit jeweiliger Rolle) anlegen" ---
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "Als Administrator kann man Mitarbeiter hinzufügen/ändern/löschen"
$ws.Range("E6").Value = ""

# --- Sprint 1 row 13: "Zimmer buchen" ---
$ws.Range("C13").Value = 1
$ws.Range("E13").Value = ""

# --- Sprint 2 rows 18-21 ---
$ws.Range("B18").Value = "Anzeige für den Koch gestalten"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 'Koch kann sich anmelden, alle offenen Bestellungen einsehen, diese aktualiesieren sich alle 5 Sekunden und diese nach zubereitung auf "gemacht" setzen '

$ws.Range("B19").Value = "Admin Fenster erweitern"
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = "Admin hat nun die Möglichkeit Speißen hinzuzufügen, bearbeiten, löschen"

$ws.Range("B20").Value = "Anzeige für Barkeeper gestalten"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 'Barkeeper kann sich anmelden, alle offen Zimmerservicebestellungen einsehen, diese aktualisieren sich alle 5 Sekunden und diese nach Zubereitung auf "fertig" setzen.'

$ws.Range("B21").Value = "Anwendung vollständig auf WebServices umstellen"
$ws.Range("C21").Value = 0.8
$ws.Range("E21").Value = "Jeglicher Datenbankzugriff soll mittels WebService funktionieren."
$ws.Range("D21").Value = "Alles wurde erfolgreich umgstellt, steckt jedoch noch in der Testphase."

# Column B grew a hair wider after the edit.
$ws.Columns.Item(2).ColumnWidth = 82.5

# Cursor/selection ends up on D21, the last cell touched.
$ws.Range("D21").Select() | Out-Null
